# IIT_Conclave4_Architecture.pptx edit:
#  1. Refresh the cached "datetimeFigureOut" footer field from 17-01-2026
#     to 18-01-2026 on the Slide Master and on every Custom Layout that
#     carries a Date placeholder.
#  2. On slide 1, rename the "Similarity Search & Filtering" caption to
#     "Semantic Search & Filtering" and shrink the auto-fit textbox to its
#     new (slightly narrower) width.

$p = $ppt.ActivePresentation

$oldDate = "17-01-2026"
$newDate = "18-01-2026"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Custom Layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide content: rename the caption textbox and re-fit its width.
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Similarity Search & Filtering") {
        $shp.TextFrame.TextRange.Text = "Semantic Search & Filtering"
        # Shape uses wrap="none" + spAutoFit, so PowerPoint re-measures the
        # box to hug the new (narrower) text; target width is 1617751 EMU.
        $shp.Width = (1617751 + 0.5) / 12700
    }
}
